$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 entirely (duplicate / bad record with malformed document number
# 113064384129). Everything below shifts up by one row.
$ws.Rows(7).Delete()

# Row 6 (Fabio Robinson Maldonado Ordoñez) now resolves: the instructor name
# is filled in and the comparison result becomes VERIFICADO.
$ws.Range("G6").Value = "FABIO ROBINSON MALDONADO ORDOÑEZ"
$ws.Range("I6").Value = "VERIFICADO"

# New row 14 (previously row 15, Gloria Yenny Castillo España) - a typo is
# introduced in the instructors name column, producing a name mismatch.
$ws.Range("G14").Value = "GLORIA YENNY CASTILLO ESPAÑAS"
$ws.Range("I14").Value = "Diferencia en nombre"

# New row 18 (previously row 19, Martha Isabel Benavides Acosta) now looks
# like it is not present in Sofía: the Ficha/Tipo/Nivel/Denominación and the
# Sofía name columns are cleared, and the comparison result is updated.
$ws.Range("A18").Value = ""
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("H18").Value = ""
$ws.Range("I18").Value = "No especificado en Sofía; Diferencia en nombre"
